$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 7257.48917153721
$ws.Range("F2").Value = -20.8986076177283
$ws.Range("C3").Value = 6939.89771781602
$ws.Range("F3").Value = 255.464211153741
$ws.Range("D5").Value = 9004
$ws.Range("E5").Value = 3297.43979530182
$ws.Range("F5").Value = -93.7576522476046
$ws.Range("D6").Value = 9004
$ws.Range("E6").Value = 3257.76184623033
$ws.Range("F6").Value = -87.1129923165712
$ws.Range("D7").Value = 9004
$ws.Range("E7").Value = 4215.66040054249
$ws.Range("F7").Value = 90.04485323912
$ws.Range("D8").Value = 9004
$ws.Range("E8").Value = 4596.2007938609
$ws.Range("F8").Value = 134.646034676944
$ws.Range("C9").Value = 7647.44774894666
$ws.Range("D9").Value = 9004
$ws.Range("E9").Value = 4596.2007938609
$ws.Range("F9").Value = 134.985355950315
$ws.Range("C10").Value = 7749.65627294543
$ws.Range("D10").Value = 9004
$ws.Range("E10").Value = 4598.15712188851
$ws.Range("F10").Value = 139.325558118081
$ws.Range("C11").Value = 7079.59727885324
$ws.Range("D11").Value = 9004
$ws.Range("E11").Value = 4628.79302058475
$ws.Range("F11").Value = 112.682929143249
$ws.Range("C12").Value = 4898.02564324643
$ws.Range("D12").Value = 9004
$ws.Range("E12").Value = 4317.45990498165
$ws.Range("F12").Value = 8.81189784283682
$ws.Range("C13").Value = 4955.91953044642
$ws.Range("D13").Value = 9004
$ws.Range("E13").Value = 4310.19379706024
$ws.Range("F13").Value = 10.921388646111
$ws.Range("C14").Value = 8643.37162142728
$ws.Range("D14").Value = 9004
$ws.Range("E14").Value = 5258.13060995135
$ws.Range("F14").Value = 204.06259297411
$ws.Range("C15").Value = 9238.96137694889
$ws.Range("D15").Value = 9004
$ws.Range("E15").Value = 5258.13437405638
$ws.Range("F15").Value = 228.878989625219
